$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.464.40'
$ws.Range("E2").Value = '  -4.64%  '

$ws.Range("D3").Value = '3.283.13'
$ws.Range("E3").Value = '  -7.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.48'
$ws.Range("E5").Value = '  -4.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '152.67'
$ws.Range("E6").Value = '  -11.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.09%  '

$ws.Range("D8").Value = '3.275.20'
$ws.Range("E8").Value = '  -7.26%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  -10.58%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.173'
$ws.Range("E10").Value = '  -13.32%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.72'
$ws.Range("E11").Value = '  -6.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.515'
$ws.Range("E12").Value = '  -12.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.75'
$ws.Range("E13").Value = '  -16.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000246'
$ws.Range("E14").Value = '  -10.74%  '

$ws.Range("D15").Value = '3.805.15'
$ws.Range("E15").Value = '  -7.09%  '

$ws.Range("D16").Value = '67.392.06'
$ws.Range("E16").Value = '  -4.85%  '

$ws.Range("D17").Value = '3.278.48'
$ws.Range("E17").Value = '  -7.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.30'
$ws.Range("E18").Value = '  -13.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '539.01'
$ws.Range("E19").Value = '  -11.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.20'
$ws.Range("E21").Value = '  -14.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.767'
$ws.Range("E22").Value = '  -12.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.95'
$ws.Range("E23").Value = '  -12.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '86.20'
$ws.Range("E24").Value = '  -11.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.67'
$ws.Range("E25").Value = '  -12.23%  '

$ws.Range("E26").Value = '  -0.10%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.26'
$ws.Range("E27").Value = '  -12.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.14'
$ws.Range("E28").Value = '  -10.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '29.53'
$ws.Range("E29").Value = '  -12.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.16'
$ws.Range("E30").Value = '  -15.54%  '

$ws.Range("E31").Value = '  -10.64%  '

$ws.Range("E32").Value = '  -11.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '546.71'
$ws.Range("E33").Value = '  -11.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.67'
$ws.Range("E34").Value = '  -17.91%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.80'
$ws.Range("E35").Value = '  -14.81%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0468'
$ws.Range("E36").Value = '  -4.95%  '

$ws.Range("B37").Value = 'FirstDigitalUSD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '53.67'
$ws.Range("E38").Value = '  -5.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0867'
$ws.Range("E39").Value = '  -12.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.16'
$ws.Range("E40").Value = '  -15.63%  '

$ws.Range("E41").Value = '  -9.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.76'
$ws.Range("E42").Value = '  -18.67%  '

$ws.Range("D43").Value = '2.948.98'
$ws.Range("E43").Value = '  -11.79%  '

$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.269'
$ws.Range("E44").Value = '  -13.15%  '

$ws.Range("B45").Value = 'PEPE'
$ws.Range("C45").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D45").Value = '0.0₃0597'
$ws.Range("E45").Value = '  -17.41%  '

$ws.Range("E46").Value = '  -15.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.20'
$ws.Range("E47").Value = '  -11.97%  '

$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.36'
$ws.Range("E49").Value = '  -18.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '126.22'
$ws.Range("E50").Value = '  -5.73%  '

$ws.Range("E51").Value = '  -11.98%  '

Write-Host "Applied cryptos update"